# fix algoritm with exception company by words from Config.xslx
#
# Insert a new settings row ("RemoveCompany") right after the existing
# "Remove_IP_UP" row, update two config values (the "prof" search mapping
# and the duplicated resume-search URL), rename the duplicate resumeURL
# entry, and add a True/False dropdown (data validation) to B7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 8 (pushes everything from the old row 8 onward down by one).
$ws.Rows(8).Insert()

# New row 8: RemoveCompany setting.
$ws.Range("A8").Value = "RemoveCompany"
$ws.Range("B8").Value = "поликлиника;больница;спартак;рапа"
$ws.Range("C8").Value = "Исключить компании, содержащие слова"

# Row 10 (was row 9 before the insert): update the "prof" mapping value.
$ws.Range("B10").Value = "Продавец,B23,urlForSalesAss"

# Row 16 (was row 15): fix the resume search URL (drop "_or_relocation").
$ws.Range("B16").Value = "https://rabota.by/search/resume?area=1002&clusters=true&currency_code=BYR&exp_period=all_time&items_on_page=100&label=only_with_salary&logic=normal&no_magichttps://rabota.by/search/resume?area=1002&label=only_with_salary&relocation=living&age_to=60&gender=unknown&clusters=true&exp_period=all_time&items_on_page=100&logic=normal&no_magic=true&order_by=relevance&ored_clusters=true&search_period=30&pos=position&text=!"

# Row 18 (was row 17): rename the duplicate "resumeURL" key so it's distinct.
$ws.Range("A18").Value = "resumeURL_copy"

# Add a True/False dropdown list validation to B7 (Remove_IP_UP).
$ws.Range("B7").Validation.Add(3, 1, 1, '"True,False"')

# Leave the active cell on B9, matching the edited file's cached selection.
$ws.Range("B9").Select()
